# Apply a cyclic rotation of data among rows 2, 3 and 4:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3
# Only columns A, B, E, F, G, H, Q, R actually change value (other columns
# happen to be identical across these three rows already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the current ("before") values for the columns of interest in rows 2-4.
$before = @{}
foreach ($r in 2..4) {
    $before[$r] = @{}
    foreach ($col in $cols) {
        $before[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Target mapping: new row -> source row (cyclic shift)
$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $before[$srcRow][$col]
    }
}
